# Apply updated crypto price/volume data (GitHub Actions symbol-list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" '304.86'
Set-TextValue "E2" '-0.18%'
Set-TextValue "D3" '35.85'
Set-TextValue "E3" '-1.29%'
Set-TextValue "D4" '4.978'
Set-TextValue "E4" '-1.94%'
Set-TextValue "D5" '0.08081'
Set-TextValue "E5" '-1.46%'
Set-TextValue "D6" '1.902'
Set-TextValue "E6" '-2.38%'
Set-TextValue "D7" '4.144'
Set-TextValue "E7" '1.78%'
Set-TextValue "D8" '7.881'
Set-TextValue "E8" '0.47%'
Set-TextValue "E9" '-0.36%'
Set-TextValue "D10" '0.1270'
Set-TextValue "E10" '-17.11%'
Set-TextValue "D11" '0.1903'
Set-TextValue "E11" '-0.93%'
Set-TextValue "D12" '0.09179'
Set-TextValue "E12" '1.29%'
Set-TextValue "D13" '0.03508'
Set-TextValue "E13" '1.91%'
Set-TextValue "D14" '0.09922'
Set-TextValue "E14" '0.65%'
Set-TextValue "D15" '0.001419'
Set-TextValue "E15" '-1.12%'
Set-TextValue "D16" '0.006069'
Set-TextValue "E16" '5.10%'
Set-TextValue "D17" '3.604'
Set-TextValue "E17" '1.82%'
Set-TextValue "D19" '0.3452'
Set-TextValue "E19" '0.19%'
Set-TextValue "D20" '5.235'
Set-TextValue "E20" '4.31%'
Set-TextValue "D21" '0.1295'
Set-TextValue "E21" '0.91%'
Set-TextValue "D22" '0.2529'
Set-TextValue "E22" '6.19%'
Set-TextValue "D23" '0.04409'
Set-TextValue "E23" '-1.54%'
Set-TextValue "D24" '0.001235'
Set-TextValue "E24" '2.76%'
Set-TextValue "D25" '0.004717'
Set-TextValue "E25" '-3.15%'
Set-TextValue "D26" '0.0001300'
Set-TextValue "E26" '6.45%'
Set-TextValue "D27" '0.0003127'
Set-TextValue "E27" '-29.02%'
Set-TextValue "D39" '0.01957'
Set-TextValue "E39" '-2.38%'
Set-TextValue "D40" '0.05226'
Set-TextValue "E40" '7.94%'
Set-TextValue "D41" '0.007553'
Set-TextValue "E41" '1.75%'
Set-TextValue "D42" '0.01015'
Set-TextValue "E42" '-2.24%'
Set-TextValue "D43" '0.1370'
Set-TextValue "E43" '0.46%'
Set-TextValue "D44" '0.002101'
Set-TextValue "E44" '0.23%'
Set-TextValue "D45" '0.01067'
Set-TextValue "E45" '-0.20%'
Set-TextValue "D46" '0.00006353'
Set-TextValue "E46" '4.77%'
Set-TextValue "D47" '0.00000000750'
Set-TextValue "E47" '0.76%'
Set-TextValue "D48" '64.96'
Set-TextValue "E48" '0.45%'
Set-TextValue "D49" '0.001658'
Set-TextValue "E49" '40.41%'
Set-TextValue "D50" '0.00002099'
Set-TextValue "E50" '0.76%'
Set-TextValue "D51" '0.0001999'
Set-TextValue "E51" '0.76%'
